$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns D (Fecha, date serial), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), S (Precio $/Kg) for rows 2..19.
$data = @{
    2  = @{ D = 44490; M = 160; N = 11500; O = 12000; P = 11750; S = 5875 }
    3  = @{ D = 44459; M = 160; N = 13000; O = 14000; P = 13500; S = 6750 }
    4  = @{ D = 44466; M = 160; N = 13500; O = 14000; P = 13750; S = 6875 }
    5  = @{ D = 44455; M = 160; N = 13000; O = 14000; P = 13500; S = 6750 }
    6  = @{ D = 44489; M = 400; N = 11500; O = 12000; P = 11750; S = 5875 }
    7  = @{ D = 44497; M = 400; N = 11500; O = 12000; P = 11750; S = 5875 }
    8  = @{ D = 44498; M = 240; N = 11000; O = 11500; P = 11250; S = 5625 }
    9  = @{ D = 44452; M = 200; N = 13000; O = 14000; P = 13500; S = 6750 }
    10 = @{ D = 44454; M = 300; N = 13000; O = 14000; P = 13500; S = 6750 }
    11 = @{ D = 44494; M = 200; N = 11500; O = 12000; P = 11750; S = 5875 }
    12 = @{ D = 44446; M = 300; N = 14000; O = 15000; P = 14500; S = 7250 }
    13 = @{ D = 44463; M = 100; N = 13000; O = 14000; P = 13500; S = 6750 }
    14 = @{ D = 44445; M = 160; N = 14000; O = 15000; P = 14500; S = 7250 }
    15 = @{ D = 44462; M = 140; N = 13000; O = 14000; P = 13500; S = 6750 }
    16 = @{ D = 44448; M = 100; N = 14000; O = 15000; P = 14500; S = 7250 }
    17 = @{ D = 44491; M = 200; N = 11500; O = 12000; P = 11750; S = 5875 }
    18 = @{ D = 44495; M = 300; N = 11000; O = 12000; P = 11500; S = 5750 }
    19 = @{ D = 44468; M = 300; N = 13000; O = 14000; P = 13500; S = 6750 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value2 = $vals.D
    $ws.Range("M$row").Value2 = $vals.M
    $ws.Range("N$row").Value2 = $vals.N
    $ws.Range("O$row").Value2 = $vals.O
    $ws.Range("P$row").Value2 = $vals.P
    $ws.Range("S$row").Value2 = $vals.S
}
